$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "58.944.20"
$ws.Cells.Item(2, 5).Value = "  +7.70%  "
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "2.527.17"
$ws.Cells.Item(3, 5).Value = "  +8.19%  "
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.00"
$ws.Cells.Item(4, 5).Value = "  +0.17%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "502.97"
$ws.Cells.Item(5, 5).Value = "  +6.81%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "155.72"
$ws.Cells.Item(6, 5).Value = "  +8.74%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.621"
$ws.Cells.Item(7, 5).Value = "  +17.26%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.992"
$ws.Cells.Item(8, 5).Value = "  -0.83%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "2.568.59"
$ws.Cells.Item(9, 5).Value = "  +9.88%  "
$ws.Cells.Item(10, 5).Value = "  +14.07%  "
$ws.Cells.Item(11, 5).Value = "  +6.98%  "
$ws.Cells.Item(12, 5).Value = "  +6.50%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.126"
$ws.Cells.Item(13, 5).Value = "  +1.66%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "2.972.74"
$ws.Cells.Item(14, 5).Value = "  +7.97%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "58.911.50"
$ws.Cells.Item(15, 5).Value = "  +7.50%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "21.77"
$ws.Cells.Item(16, 5).Value = "  +9.14%  "
$ws.Cells.Item(17, 5).Value = "  +5.61%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "2.646.54"
$ws.Cells.Item(18, 5).Value = "  +13.06%  "
$ws.Cells.Item(19, 5).Value = "  +4.43%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "332.38"
$ws.Cells.Item(20, 5).Value = "  +6.50%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "10.30"
$ws.Cells.Item(21, 5).Value = "  +7.71%  "
$ws.Cells.Item(22, 5).Value = "  +7.82%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "0.996"
$ws.Cells.Item(23, 5).Value = "  -0.14%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "59.44"
$ws.Cells.Item(24, 5).Value = "  +5.74%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "0.415"
$ws.Cells.Item(25, 5).Value = "  +6.44%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "0.165"
$ws.Cells.Item(26, 5).Value = "  +8.01%  "
$ws.Cells.Item(27, 2).Value = "WrappedeETH"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "2.662.72"
$ws.Cells.Item(27, 5).Value = "  +8.80%  "
$ws.Cells.Item(28, 2).Value = "Binance-PegBSC-USD"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "0.995"
$ws.Cells.Item(28, 5).Value = "  -0.40%  "
$ws.Cells.Item(29, 2).Value = "PEPE"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "0.0₃0823"
$ws.Cells.Item(29, 5).Value = "  +10.26%  "
$ws.Cells.Item(30, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "7.39"
$ws.Cells.Item(30, 5).Value = "  +3.95%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "0.997"
$ws.Cells.Item(31, 5).Value = "  -0.34%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "157.27"
$ws.Cells.Item(32, 5).Value = "  +7.44%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "19.31"
$ws.Cells.Item(33, 5).Value = "  +7.08%  "
$ws.Cells.Item(34, 5).Value = "  +7.10%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "5.49"
$ws.Cells.Item(35, 5).Value = "  +9.54%  "
$ws.Cells.Item(36, 5).Value = "  +9.90%  "
$ws.Cells.Item(37, 5).Value = "  +9.26%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.841"
$ws.Cells.Item(38, 5).Value = "  +3.42%  "
$ws.Cells.Item(39, 5).Value = "  +11.81%  "
$ws.Cells.Item(40, 5).Value = "  +8.61%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "35.07"
$ws.Cells.Item(41, 5).Value = "  +5.16%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "289.89"
$ws.Cells.Item(42, 5).Value = "  +16.14%  "
$ws.Cells.Item(43, 5).Value = "  +6.96%  "
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.623"
$ws.Cells.Item(44, 5).Value = "  +8.60%  "
$ws.Cells.Item(45, 2).Value = "SuiNetwork"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.780"
$ws.Cells.Item(45, 5).Value = "  +25.23%  "
$ws.Cells.Item(46, 2).Value = "Hedera"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.0564"
$ws.Cells.Item(46, 5).Value = "  +8.28%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.990"
$ws.Cells.Item(47, 5).Value = "  -1.06%  "
$ws.Cells.Item(48, 2).Value = "RenderToken"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "4.85"
$ws.Cells.Item(48, 5).Value = "  +11.24%  "
$ws.Cells.Item(49, 2).Value = "EnergySwap"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "19.10"
$ws.Cells.Item(49, 5).Value = "  +14.61%  "
$ws.Cells.Item(50, 2).Value = "VeChain"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.0235"
$ws.Cells.Item(50, 5).Value = "  +6.94%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "2.000.31"
$ws.Cells.Item(51, 5).Value = "  +12.61%  "
